{"js": "const replacements = [\n  [\"373\u00f79=\", \"335\u00f78=\"],\n  [\"948\u00f78=\", \"307\u00f77=\"],\n  [\"532\u00f75=\", \"166\u00f74=\"],\n  [\"748\u00f74=\", \"678\u00f77=\"],\n  [\"739\u00f77=\", \"801\u00f72=\"],\n  [\"568\u00f75=\", \"612\u00f76=\"],\n  [\"939\u00f76=\", \"304\u00f75=\"],\n  [\"474\u00f74=\", \"725\u00f73=\"],\n  [\"107\u00f72=\", \"404\u00f75=\"],\n  [\"781\u00f79=\", \"261\u00f72=\"],\n  [\"123\u00f72=\", \"479\u00f74=\"],\n  [\"560\u00f73=\", \"814\u00f73=\"],\n  [\"567\u00f79=\", \"788\u00f76=\"],\n  [\"964\u00f73=\", \"437\u00f77=\"],\n  [\"266\u00f77=\", \"962\u00f78=\"],\n  [\"887\u00f74=\", \"629\u00f72=\"],\n  [\"523\u00f72=\", \"595\u00f76=\"],\n  [\"216\u00f74=\", \"461\u00f76=\"],\n  [\"942\u00f72=\", \"746\u00f75=\"],\n  [\"947\u00f76=\", \"972\u00f74=\"],\n  [\"837\u00f79=\", \"995\u00f73=\"],\n  [\"924\u00f73=\", \"140\u00f77=\"],\n  [\"225\u00f74=\", \"776\u00f73=\"],\n  [\"419\u00f73=\", \"361\u00f78=\"],\n  [\"389\u00f78=\", \"878\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"373\u00f79=\", \"335\u00f78=\"),\n    @(\"948\u00f78=\", \"307\u00f77=\"),\n    @(\"532\u00f75=\", \"166\u00f74=\"),\n    @(\"748\u00f74=\", \"678\u00f77=\"),\n    @(\"739\u00f77=\", \"801\u00f72=\"),\n    @(\"568\u00f75=\", \"612\u00f76=\"),\n    @(\"939\u00f76=\", \"304\u00f75=\"),\n    @(\"474\u00f74=\", \"725\u00f73=\"),\n    @(\"107\u00f72=\", \"404\u00f75=\"),\n    @(\"781\u00f79=\", \"261\u00f72=\"),\n    @(\"123\u00f72=\", \"479\u00f74=\"),\n    @(\"560\u00f73=\", \"814\u00f73=\"),\n    @(\"567\u00f79=\", \"788\u00f76=\"),\n    @(\"964\u00f73=\", \"437\u00f77=\"),\n    @(\"266\u00f77=\", \"962\u00f78=\"),\n    @(\"887\u00f74=\", \"629\u00f72=\"),\n    @(\"523\u00f72=\", \"595\u00f76=\"),\n    @(\"216\u00f74=\", \"461\u00f76=\"),\n    @(\"942\u00f72=\", \"746\u00f75=\"),\n    @(\"947\u00f76=\", \"972\u00f74=\"),\n    @(\"837\u00f79=\", \"995\u00f73=\"),\n    @(\"924\u00f73=\", \"140\u00f77=\"),\n    @(\"225\u00f74=\", \"776\u00f73=\"),\n    @(\"419\u00f73=\", \"361\u00f78=\"),\n    @(\"389\u00f78=\", \"878\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
